$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.59%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.49%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.932"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'15.93%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08139"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.38%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.589"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.51%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.740"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.33%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.945"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.15%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D10").Value = "'0.9435"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.24%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1307"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'16.65%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'7.39%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.04%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03447"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'4.59%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09608"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.83%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001314"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.51%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006541"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'9.14%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.372"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.07%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3533"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.50%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.739"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'21.78%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1450"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'11.57%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2449"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.72%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04453"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.71%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001253"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'4.14%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004353"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.19%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001190"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-10.76%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003988"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'37.23%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02513"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'19.43%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05328"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'8.71%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007606"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.05%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1432"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.22%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008942"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.46%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002065"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.46%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'26.63%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006777"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'3.89%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.25%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002895"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-12.35%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001799"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'24.43%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.25%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.25%"
$ws.Range("E51").Style = "Normal"
